# Revert "Merge branch 'main' into waowaowiwaow"
# -------------------------------------------------------------
# Substantive data changes being restored:
#   F2 (SURNAME)             : "Vinoya"             -> "Vinoya "
#   N2 (BENEFICIARY FULL NAME): "Ann Michel Pascual" -> "Ann Michel Pascual "
# Plus the cursor/selection position that was active when the
# earlier (reverted) revision was saved: cell Q9 (with the view
# scrolled so column M is the left-most visible column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- restore the trailing-space surname/beneficiary-name values ---
$ws.Range("F2").Value = "Vinoya "
$ws.Range("N2").Value = "Ann Michel Pascual "

# --- restore the view/selection state (topLeftCell=M1, activeCell/sqref=Q9) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 13

$ws.Range("Q9").Select()
